$d = $word.ActiveDocument

# Locate the paragraph ending in "Client disconnect -> ket thuc stream" via Find
# (ASCII-safe substring so the match is robust regardless of host text encoding).
$findRange = $d.Content
$found = $findRange.Find.Execute("Client disconnect", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchPos = $findRange.Start

# Resolve the reliable, document-wide 1-based paragraph index for the match
# (Range.Paragraphs.Item/.First is not trustworthy in this host, so scan instead).
$targetIndex = -1
$scanIdx = 0
foreach ($p in $d.Paragraphs) {
    $scanIdx = $scanIdx + 1
    if ($p.Range.Start -le $matchPos -and $matchPos -lt $p.Range.End) {
        $targetIndex = $scanIdx
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$insertAt = $target.Range.End

# --- Insert a new, blank paragraph right after it ---------------------------------
$blankPoint = $d.Range($insertAt, $insertAt)
$blankPoint.InsertParagraphBefore()

# --- Insert the new "Database Design" Heading2 paragraph after the blank one ------
$blankPara = $d.Paragraphs.Item($targetIndex + 1)
$blankEnd = $blankPara.Range.End
$headingPoint = $d.Range($blankEnd, $blankEnd)
$headingPoint.InsertParagraphBefore()

$headingPara = $d.Paragraphs.Item($targetIndex + 2)
$headingPara.Style = "Heading2"

$headingStart = $headingPara.Range.Start
$headingTextPoint = $d.Range($headingStart, $headingStart)
$headingTextPoint.Text = "Database"

$secondRunAt = $headingStart + 8
$secondRunPoint = $d.Range($secondRunAt, $secondRunAt)
$secondRunPoint.InsertAfter(" Design")
